$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Col4a6"
$ws.Cells.Item(2,3).Value = "Cd93"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.1072476666666667
$ws.Cells.Item(2,8).Value = 0.321743
$ws.Cells.Item(2,9).Value = 0.0853607675234147
$ws.Cells.Item(2,10).Value = 0.08536076752341469
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 209.26237
$ws.Cells.Item(2,14).Value = 627.78711
$ws.Cells.Item(2,15).Value = 0.8127157202241573
$ws.Cells.Item(2,16).Value = 0.8127157202241573
$ws.Cells.Item(2,17).Value = 22.44290090363667
$ws.Cells.Item(2,18).Value = 201.98610813273
$ws.Cells.Item(2,19).Value = 0.06937403765667884
$ws.Cells.Item(2,20).Value = 0.06937403765667882

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Col4a6"
$ws.Cells.Item(3,3).Value = "Cd93"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.1072476666666667
$ws.Cells.Item(3,8).Value = 0.321743
$ws.Cells.Item(3,9).Value = 0.0853607675234147
$ws.Cells.Item(3,10).Value = 0.08536076752341469
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.9848756666666668
$ws.Cells.Item(3,14).Value = 2.954627
$ws.Cells.Item(3,15).Value = 0.003824977881910862
$ws.Cells.Item(3,16).Value = 0.003824977881910862
$ws.Cells.Item(3,17).Value = 0.1056256172067778
$ws.Cells.Item(3,18).Value = 0.9506305548610001
$ws.Cells.Item(3,19).Value = 0.0003265030477599963
$ws.Cells.Item(3,20).Value = 0.0003265030477599962

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Col4a6"
$ws.Cells.Item(4,3).Value = "Cd93"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.1072476666666667
$ws.Cells.Item(4,8).Value = 0.321743
$ws.Cells.Item(4,9).Value = 0.0853607675234147
$ws.Cells.Item(4,10).Value = 0.08536076752341469
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.763846666666667
$ws.Cells.Item(4,14).Value = 5.291539999999999
$ws.Cells.Item(4,15).Value = 0.006850280411451801
$ws.Cells.Item(4,16).Value = 0.006850280411451801
$ws.Cells.Item(4,17).Value = 0.1891684393577778
$ws.Cells.Item(4,18).Value = 1.70251595422
$ws.Cells.Item(4,19).Value = 0.0005847451936721388
$ws.Cells.Item(4,20).Value = 0.0005847451936721387

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Col4a6"
$ws.Cells.Item(5,3).Value = "Cd93"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.1072476666666667
$ws.Cells.Item(5,8).Value = 0.321743
$ws.Cells.Item(5,9).Value = 0.0853607675234147
$ws.Cells.Item(5,10).Value = 0.08536076752341469
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 45.474231
$ws.Cells.Item(5,14).Value = 136.422693
$ws.Cells.Item(5,15).Value = 0.1766090214824801
$ws.Cells.Item(5,16).Value = 0.1766090214824801
$ws.Cells.Item(5,17).Value = 4.877005168210999
$ws.Cells.Item(5,18).Value = 43.89304651389899
$ws.Cells.Item(5,19).Value = 0.01507548162530374
$ws.Cells.Item(5,20).Value = 0.01507548162530373

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Col4a6"
$ws.Cells.Item(6,3).Value = "Cd93"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.3710633333333333
$ws.Cells.Item(6,8).Value = 1.11319
$ws.Cells.Item(6,9).Value = 0.2953374363992068
$ws.Cells.Item(6,10).Value = 0.2953374363992068
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 209.26237
$ws.Cells.Item(6,14).Value = 627.78711
$ws.Cells.Item(6,15).Value = 0.8127157202241573
$ws.Cells.Item(6,16).Value = 0.8127157202241573
$ws.Cells.Item(6,17).Value = 77.64959255343332
$ws.Cells.Item(6,18).Value = 698.8463329808999
$ws.Cells.Item(6,19).Value = 0.2400253773323376
$ws.Cells.Item(6,20).Value = 0.2400253773323376

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Col4a6"
$ws.Cells.Item(7,3).Value = "Cd93"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.3710633333333333
$ws.Cells.Item(7,8).Value = 1.11319
$ws.Cells.Item(7,9).Value = 0.2953374363992068
$ws.Cells.Item(7,10).Value = 0.2953374363992068
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.9848756666666668
$ws.Cells.Item(7,14).Value = 2.954627
$ws.Cells.Item(7,15).Value = 0.003824977881910862
$ws.Cells.Item(7,16).Value = 0.003824977881910862
$ws.Cells.Item(7,17).Value = 0.3654512477922223
$ws.Cells.Item(7,18).Value = 3.28906123013
$ws.Cells.Item(7,19).Value = 0.001129659161927222
$ws.Cells.Item(7,20).Value = 0.001129659161927222

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Col4a6"
$ws.Cells.Item(8,3).Value = "Cd93"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.3710633333333333
$ws.Cells.Item(8,8).Value = 1.11319
$ws.Cells.Item(8,9).Value = 0.2953374363992068
$ws.Cells.Item(8,10).Value = 0.2953374363992068
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 1.763846666666667
$ws.Cells.Item(8,14).Value = 5.291539999999999
$ws.Cells.Item(8,15).Value = 0.006850280411451801
$ws.Cells.Item(8,16).Value = 0.006850280411451801
$ws.Cells.Item(8,17).Value = 0.6544988236222221
$ws.Cells.Item(8,18).Value = 5.890489412599999
$ws.Cells.Item(8,19).Value = 0.002023144255333879
$ws.Cells.Item(8,20).Value = 0.002023144255333878

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Col4a6"
$ws.Cells.Item(9,3).Value = "Cd93"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.3710633333333333
$ws.Cells.Item(9,8).Value = 1.11319
$ws.Cells.Item(9,9).Value = 0.2953374363992068
$ws.Cells.Item(9,10).Value = 0.2953374363992068
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 45.474231
$ws.Cells.Item(9,14).Value = 136.422693
$ws.Cells.Item(9,15).Value = 0.1766090214824801
$ws.Cells.Item(9,16).Value = 0.1766090214824801
$ws.Cells.Item(9,17).Value = 16.87381973563
$ws.Cells.Item(9,18).Value = 151.86437762067
$ws.Cells.Item(9,19).Value = 0.05215925564960812
$ws.Cells.Item(9,20).Value = 0.05215925564960811

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Col4a6"
$ws.Cells.Item(10,3).Value = "Cd93"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.7657353333333333
$ws.Cells.Item(10,8).Value = 2.297206
$ws.Cells.Item(10,9).Value = 0.6094655278262259
$ws.Cells.Item(10,10).Value = 0.6094655278262258
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 209.26237
$ws.Cells.Item(10,14).Value = 627.78711
$ws.Cells.Item(10,15).Value = 0.8127157202241573
$ws.Cells.Item(10,16).Value = 0.8127157202241573
$ws.Cells.Item(10,17).Value = 160.2395906460733
$ws.Cells.Item(10,18).Value = 1442.15631581466
$ws.Cells.Item(10,19).Value = 0.4953222153990873
$ws.Cells.Item(10,20).Value = 0.4953222153990872

# Row 11
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Col4a6"
$ws.Cells.Item(11,3).Value = "Cd93"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.7657353333333333
$ws.Cells.Item(11,8).Value = 2.297206
$ws.Cells.Item(11,9).Value = 0.6094655278262259
$ws.Cells.Item(11,10).Value = 0.6094655278262258
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.9848756666666668
$ws.Cells.Item(11,14).Value = 2.954627
$ws.Cells.Item(11,15).Value = 0.003824977881910862
$ws.Cells.Item(11,16).Value = 0.003824977881910862
$ws.Cells.Item(11,17).Value = 0.754154096906889
$ws.Cells.Item(11,18).Value = 6.787386872162001
$ws.Cells.Item(11,19).Value = 0.002331192163722443
$ws.Cells.Item(11,20).Value = 0.002331192163722443

# Row 12
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Col4a6"
$ws.Cells.Item(12,3).Value = "Cd93"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.7657353333333333
$ws.Cells.Item(12,8).Value = 2.297206
$ws.Cells.Item(12,9).Value = 0.6094655278262259
$ws.Cells.Item(12,10).Value = 0.6094655278262258
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 1.763846666666667
$ws.Cells.Item(12,14).Value = 5.291539999999999
$ws.Cells.Item(12,15).Value = 0.006850280411451801
$ws.Cells.Item(12,16).Value = 0.006850280411451801
$ws.Cells.Item(12,17).Value = 1.350639715248889
$ws.Cells.Item(12,18).Value = 12.15575743724
$ws.Cells.Item(12,19).Value = 0.004175009766723127
$ws.Cells.Item(12,20).Value = 0.004175009766723127

# Row 13
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Col4a6"
$ws.Cells.Item(13,3).Value = "Cd93"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.7657353333333333
$ws.Cells.Item(13,8).Value = 2.297206
$ws.Cells.Item(13,9).Value = 0.6094655278262259
$ws.Cells.Item(13,10).Value = 0.6094655278262258
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 45.474231
$ws.Cells.Item(13,14).Value = 136.422693
$ws.Cells.Item(13,15).Value = 0.1766090214824801
$ws.Cells.Item(13,16).Value = 0.1766090214824801
$ws.Cells.Item(13,17).Value = 34.82122543286199
$ws.Cells.Item(13,18).Value = 313.391028895758
$ws.Cells.Item(13,19).Value = 0.107637110496693
$ws.Cells.Item(13,20).Value = 0.107637110496693

# Row 14
$ws.Cells.Item(14,1).Value = "Resolving-Mac"
$ws.Cells.Item(14,2).Value = "Col4a6"
$ws.Cells.Item(14,3).Value = "Cd93"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 0.3333333333333333
$ws.Cells.Item(14,7).Value = 0.01235833333333333
$ws.Cells.Item(14,8).Value = 0.037075
$ws.Cells.Item(14,9).Value = 0.009836268251152627
$ws.Cells.Item(14,10).Value = 0.009836268251152625
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 209.26237
$ws.Cells.Item(14,14).Value = 627.78711
$ws.Cells.Item(14,15).Value = 0.8127157202241573
$ws.Cells.Item(14,16).Value = 0.8127157202241573
$ws.Cells.Item(14,17).Value = 2.586134122583333
$ws.Cells.Item(14,18).Value = 23.27520710325
$ws.Cells.Item(14,19).Value = 0.00799408983605352
$ws.Cells.Item(14,20).Value = 0.007994089836053517

# Row 15
$ws.Cells.Item(15,1).Value = "Resolving-Mac"
$ws.Cells.Item(15,2).Value = "Col4a6"
$ws.Cells.Item(15,3).Value = "Cd93"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = 0.3333333333333333
$ws.Cells.Item(15,7).Value = 0.01235833333333333
$ws.Cells.Item(15,8).Value = 0.037075
$ws.Cells.Item(15,9).Value = 0.009836268251152627
$ws.Cells.Item(15,10).Value = 0.009836268251152625
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.9848756666666668
$ws.Cells.Item(15,14).Value = 2.954627
$ws.Cells.Item(15,15).Value = 0.003824977881910862
$ws.Cells.Item(15,16).Value = 0.003824977881910862
$ws.Cells.Item(15,17).Value = 0.01217142178055556
$ws.Cells.Item(15,18).Value = 0.109542796025
$ws.Cells.Item(15,19).Value = 0.00003762350850120084
$ws.Cells.Item(15,20).Value = 0.00003762350850120083

# Row 16
$ws.Cells.Item(16,1).Value = "Resolving-Mac"
$ws.Cells.Item(16,2).Value = "Col4a6"
$ws.Cells.Item(16,3).Value = "Cd93"
$ws.Cells.Item(16,4).Value = "MuSCs"
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = 0.3333333333333333
$ws.Cells.Item(16,7).Value = 0.01235833333333333
$ws.Cells.Item(16,8).Value = 0.037075
$ws.Cells.Item(16,9).Value = 0.009836268251152627
$ws.Cells.Item(16,10).Value = 0.009836268251152625
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 1.763846666666667
$ws.Cells.Item(16,14).Value = 5.291539999999999
$ws.Cells.Item(16,15).Value = 0.006850280411451801
$ws.Cells.Item(16,16).Value = 0.006850280411451801
$ws.Cells.Item(16,17).Value = 0.02179820505555555
$ws.Cells.Item(16,18).Value = 0.1961838455
$ws.Cells.Item(16,19).Value = 0.0000673811957226561
$ws.Cells.Item(16,20).Value = 0.0000673811957226561

# Row 17
$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Col4a6"
$ws.Cells.Item(17,3).Value = "Cd93"
$ws.Cells.Item(17,4).Value = "Resolving-Mac"
$ws.Cells.Item(17,5).Value = 1
$ws.Cells.Item(17,6).Value = 0.3333333333333333
$ws.Cells.Item(17,7).Value = 0.01235833333333333
$ws.Cells.Item(17,8).Value = 0.037075
$ws.Cells.Item(17,9).Value = 0.009836268251152627
$ws.Cells.Item(17,10).Value = 0.009836268251152625
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 45.474231
$ws.Cells.Item(17,14).Value = 136.422693
$ws.Cells.Item(17,15).Value = 0.1766090214824801
$ws.Cells.Item(17,16).Value = 0.1766090214824801
$ws.Cells.Item(17,17).Value = 0.5619857047749999
$ws.Cells.Item(17,18).Value = 5.057871342974999
$ws.Cells.Item(17,19).Value = 0.001737173710875251
$ws.Cells.Item(17,20).Value = 0.001737173710875251

